$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "L01"
$ws.Range("D3").Value = "L02"
$ws.Range("D4").Value = "L03"
$ws.Range("E2").Value = "RAG BASED MODEL"
$ws.Range("E3").Value = "ADVANCED LLM"
$ws.Range("E4").Value = "ADVANCED NLP"

$ws.Columns.Item(5).ColumnWidth = 22

$ws.Range("E5").Select()
